$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.254.88"
$ws.Range("E2").Value = "  -1.87%  "
# Row 3
$ws.Range("D3").Value = "2.918.61"
$ws.Range("E3").Value = "  -0.59%  "
# Row 4
$ws.Range("E4").Value = "  -0.21%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "372.16"
$ws.Range("E5").Value = "  +4.47%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.59"
$ws.Range("E6").Value = "  -5.66%  "
# Row 7
$ws.Range("E7").Value = "  -5.09%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.07%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -5.60%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.12"
$ws.Range("E10").Value = "  -4.76%  "
# Row 11
$ws.Range("E11").Value = "  +0.91%  "
# Row 12
$ws.Range("E12").Value = "  -3.38%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.43"
$ws.Range("E13").Value = "  -5.52%  "
# Row 14
$ws.Range("D14").Value = "3.380.00"
$ws.Range("E14").Value = "  -0.63%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.39"
$ws.Range("E15").Value = "  -5.13%  "
# Row 16
$ws.Range("D16").Value = "2.912.20"
$ws.Range("E16").Value = "  -0.70%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.954"
$ws.Range("E17").Value = "  -2.50%  "
# Row 18
$ws.Range("D18").Value = "51.239.52"
$ws.Range("E18").Value = "  -1.81%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.33"
$ws.Range("E19").Value = "  -6.24%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  -3.55%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("E21").Value = "  -5.50%  "
# Row 22
$ws.Range("E22").Value = "  -3.26%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.41"
$ws.Range("E23").Value = "  -2.81%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.44"
$ws.Range("E24").Value = "  -3.24%  "
# Row 25
$ws.Range("E25").Value = "  -2.40%  "
# Row 26
$ws.Range("E26").Value = "  +4.29%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.175"
$ws.Range("E27").Value = "  -2.22%  "
# Row 28
$ws.Range("E28").Value = "  +0.02%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.44"
$ws.Range("E29").Value = "  -5.71%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.96"
$ws.Range("E30").Value = "  -3.66%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.104"
$ws.Range("E31").Value = "  -3.18%  "
# Row 32
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.19"
$ws.Range("E32").Value = "  -0.21%  "
# Row 33
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.96"
$ws.Range("E33").Value = "  -4.75%  "
# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.29"
$ws.Range("E34").Value = "  -6.05%  "
# Row 35
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.14"
$ws.Range("E35").Value = "  -6.10%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.65"
$ws.Range("E36").Value = "  -2.64%  "
# Row 37
$ws.Range("E37").Value = "  +0.38%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0425"
$ws.Range("E38").Value = "  -3.96%  "
# Row 39
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  -1.58%  "
# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("E40").Value = "  +0.83%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.02"
$ws.Range("E41").Value = "  -6.28%  "
# Row 42
$ws.Range("E42").Value = "  -6.01%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.114"
$ws.Range("E43").Value = "  -5.25%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.25"
$ws.Range("E44").Value = "  -2.75%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "117.38"
$ws.Range("E45").Value = "  -2.05%  "
# Row 46
$ws.Range("E46").Value = "  -3.63%  "
# Row 47
$ws.Range("D47").Value = "2.060.12"
$ws.Range("E47").Value = "  -3.42%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -5.82%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.20"
$ws.Range("E49").Value = "  -7.22%  "
# Row 50
$ws.Range("D50").Value = "3.214.21"
$ws.Range("E50").Value = "  -0.48%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.236"
$ws.Range("E51").Value = "  -5.26%  "
